$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the target paragraph unambiguously via a stable anchor.
# ------------------------------------------------------------------
$content = $d.Content.Text
$paraAnchor = $content.IndexOf("Beispiel arithmetischer")

# ------------------------------------------------------------------
# 1) First occurrence: "Linksshift" -> "Rechtsshift"
#    (keep "Rechts" and "shift" as two separate runs, like the source)
# ------------------------------------------------------------------
$content = $d.Content.Text
$idx = $content.IndexOf("Linksshift", $paraAnchor)
$rLinks = $d.Range($idx, $idx + 5)
$rLinks.Text = "Rechts"

# Re-split "Rechts" from "shift" (the assignment above merges them into
# one run); toggling Bold off->on->off on the trailing part forces the
# serializer to keep it as its own run without leaving any stray
# formatting behind.
$content = $d.Content.Text
$idx = $content.IndexOf("Rechtsshift", $paraAnchor)
$shiftStart = $idx + 6
$shiftEnd = $shiftStart + 5
$rShift = $d.Range($shiftStart, $shiftEnd)
$rShift.Bold = $true
$rShift.Bold = $false

# ------------------------------------------------------------------
# 2) "<<<" -> ">>>"
#    result must be split into "  „" / ">>>" / "“ und logischer "
# ------------------------------------------------------------------
$content = $d.Content.Text
$idx = $content.IndexOf("<<<", $paraAnchor)
$rArrows = $d.Range($idx, $idx + 3)
$rArrows.Text = ">>>"

$content = $d.Content.Text
$idx = $content.IndexOf(">>>", $paraAnchor)
$rArrows2 = $d.Range($idx, $idx + 3)
$rArrows2.Bold = $true
$rArrows2.Bold = $false

# ------------------------------------------------------------------
# 3) Second occurrence: "Linkssshift" -> "Rechtsshift"
#    (the two runs naturally coalesce into a single run, matching target)
# ------------------------------------------------------------------
$content = $d.Content.Text
$idx = $content.IndexOf("Linkssshift", $paraAnchor)
$rLinks2 = $d.Range($idx, $idx + 11)
$rLinks2.Text = "Rechtsshift"

# ------------------------------------------------------------------
# 4) Last "<<" -> ">>" (the one right after the second "Rechtsshift")
#    result must be split into " „>>" / "“."
# ------------------------------------------------------------------
$content = $d.Content.Text
$secondRechts = $content.IndexOf("Rechtsshift", $idx + 1)
$idx = $content.IndexOf("<<", $secondRechts)
$rArrows3 = $d.Range($idx, $idx + 2)
$rArrows3.Text = ">>"

$content = $d.Content.Text
$idx = $content.IndexOf("“.", $secondRechts)
$rEnd = $d.Range($idx, $idx + 2)
$rEnd.Bold = $true
$rEnd.Bold = $false
